$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row after "member-identification.html" (row 46) so that a
# new "member-identification-log.html" entry becomes row 47, pushing all
# subsequent rows down by one.
$ws.Rows("47:47").Insert()
$ws.Cells.Item(47, 2).Value = "鑑定紀錄"
$ws.Cells.Item(47, 3).Value = "member-identification-log.html"

# Append two new rows at the bottom of the list (rows 62 and 63).
$ws.Cells.Item(62, 2).Value = "購買刊登筆數"
$ws.Cells.Item(62, 3).Value = "seller-service.html"

$ws.Cells.Item(63, 2).Value = "購買 VIP 會員方案"
$ws.Cells.Item(63, 3).Value = "seller-service-checkout.html"

# Match the final view/selection state recorded in the workbook.
$ws.Range("D63").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 49
